$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(307, 1).Value = 5
$ws.Cells.Item(307, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(307, 2).Font.Italic = $true
$ws.Cells.Item(307, 3).Value = "52mm SL"
$ws.Cells.Item(307, 4).Value = 1.65
$ws.Cells.Item(307, 5).Value = 2
$ws.Cells.Item(307, 6).Value = "F1"
$ws.Cells.Item(307, 7).Value = "TC5"
$ws.Cells.Item(307, 8).Value = "TCMM214"

$ws.Cells.Item(308, 1).Value = 5
$ws.Cells.Item(308, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(308, 2).Font.Italic = $true
$ws.Cells.Item(308, 3).Value = "48mm SL"
$ws.Cells.Item(308, 4).Value = 1.18
$ws.Cells.Item(308, 5).Value = 2
$ws.Cells.Item(308, 6).Value = "F2"
$ws.Cells.Item(308, 7).Value = "TC5"
$ws.Cells.Item(308, 8).Value = "TCMM215"

$ws.Cells.Item(309, 1).Value = 5
$ws.Cells.Item(309, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(309, 2).Font.Italic = $true
$ws.Cells.Item(309, 3).Value = "49mm SL"
$ws.Cells.Item(309, 4).Value = 1.45
$ws.Cells.Item(309, 5).Value = 1
$ws.Cells.Item(309, 6).Value = "F3"
$ws.Cells.Item(309, 7).Value = "TC5"
$ws.Cells.Item(309, 8).Value = "TCMM216"

$ws.Cells.Item(310, 1).Value = 5
$ws.Cells.Item(310, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(310, 2).Font.Italic = $true
$ws.Cells.Item(310, 3).Value = "50mm SL"
$ws.Cells.Item(310, 4).Value = 1.48
$ws.Cells.Item(310, 5).Value = 2
$ws.Cells.Item(310, 6).Value = "F4"
$ws.Cells.Item(310, 7).Value = "TC5"
$ws.Cells.Item(310, 8).Value = "TCMM217"

$ws.Cells.Item(311, 1).Value = 5
$ws.Cells.Item(311, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(311, 2).Font.Italic = $true
$ws.Cells.Item(311, 3).Value = "47mm SL"
$ws.Cells.Item(311, 4).Value = 1.21
$ws.Cells.Item(311, 5).Value = 2
$ws.Cells.Item(311, 6).Value = "F5"
$ws.Cells.Item(311, 7).Value = "TC5"
$ws.Cells.Item(311, 8).Value = "TCMM218"

$ws.Cells.Item(312, 1).Value = 5
$ws.Cells.Item(312, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(312, 2).Font.Italic = $true
$ws.Cells.Item(312, 3).Value = "54mm SL"
$ws.Cells.Item(312, 4).Value = 2.21
$ws.Cells.Item(312, 5).Value = 2
$ws.Cells.Item(312, 6).Value = "F6"
$ws.Cells.Item(312, 7).Value = "TC5"
$ws.Cells.Item(312, 8).Value = "TCMM219"

$ws.Cells.Item(313, 1).Value = 5
$ws.Cells.Item(313, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(313, 2).Font.Italic = $true
$ws.Cells.Item(313, 3).Value = "50mm SL"
$ws.Cells.Item(313, 4).Value = 1.3
$ws.Cells.Item(313, 5).Value = 2
$ws.Cells.Item(313, 6).Value = "F7"
$ws.Cells.Item(313, 7).Value = "TC5"
$ws.Cells.Item(313, 8).Value = "TCMM220"

$ws.Cells.Item(314, 1).Value = 5
$ws.Cells.Item(314, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(314, 2).Font.Italic = $true
$ws.Cells.Item(314, 3).Value = "47mm SL"
$ws.Cells.Item(314, 4).Value = 1.16
$ws.Cells.Item(314, 5).Value = 2
$ws.Cells.Item(314, 6).Value = "F8"
$ws.Cells.Item(314, 7).Value = "TC5"
$ws.Cells.Item(314, 8).Value = "TCMM221"

$ws.Cells.Item(315, 1).Value = 5
$ws.Cells.Item(315, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(315, 2).Font.Italic = $true
$ws.Cells.Item(315, 3).Value = "48mm SL"
$ws.Cells.Item(315, 4).Value = 1.37
$ws.Cells.Item(315, 5).Value = 2
$ws.Cells.Item(315, 6).Value = "F9"
$ws.Cells.Item(315, 7).Value = "TC5"
$ws.Cells.Item(315, 8).Value = "TCMM222"

$ws.Cells.Item(316, 1).Value = 5
$ws.Cells.Item(316, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(316, 2).Font.Italic = $true
$ws.Cells.Item(316, 3).Value = "50mm SL"
$ws.Cells.Item(316, 4).Value = 1.53
$ws.Cells.Item(316, 5).Value = 2
$ws.Cells.Item(316, 6).Value = "F10"
$ws.Cells.Item(316, 7).Value = "TC5"
$ws.Cells.Item(316, 8).Value = "TCMM223"

$ws.Cells.Item(317, 1).Value = 5
$ws.Cells.Item(317, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(317, 2).Font.Italic = $true
$ws.Cells.Item(317, 3).Value = "50mm SL"
$ws.Cells.Item(317, 4).Value = 1.42
$ws.Cells.Item(317, 5).Value = 2
$ws.Cells.Item(317, 6).Value = "F11"
$ws.Cells.Item(317, 7).Value = "TC5"
$ws.Cells.Item(317, 8).Value = "TCMM224"

$ws.Cells.Item(318, 1).Value = 5
$ws.Cells.Item(318, 2).Value = "Maurolicus muelleri"
$ws.Cells.Item(318, 2).Font.Italic = $true
$ws.Cells.Item(318, 3).Value = "47mm SL"
$ws.Cells.Item(318, 4).Value = 1.22
$ws.Cells.Item(318, 5).Value = 2
$ws.Cells.Item(318, 6).Value = "F12"
$ws.Cells.Item(318, 7).Value = "TC5"
$ws.Cells.Item(318, 8).Value = "TCMM225"

$ws.Cells.Item(319, 1).Value = 2
$ws.Cells.Item(319, 2).Value = "Argyropelecus olfersii"
$ws.Cells.Item(319, 2).Font.Italic = $true
$ws.Cells.Item(319, 3).Value = "44mm SL"
$ws.Cells.Item(319, 4).Value = 3.03
$ws.Cells.Item(319, 5).Value = 1
$ws.Cells.Item(319, 6).Value = "B6"
$ws.Cells.Item(319, 7).Value = "TC4"
$ws.Cells.Item(319, 8).Value = "TCAO016"

$ws.Cells.Item(320, 1).Value = 2
$ws.Cells.Item(320, 2).Value = "Argyropelecus olfersii"
$ws.Cells.Item(320, 2).Font.Italic = $true
$ws.Cells.Item(320, 3).Value = "61mm SL"
$ws.Cells.Item(320, 4).Value = 5.43
$ws.Cells.Item(320, 5).Value = 1
$ws.Cells.Item(320, 6).Value = "B7"
$ws.Cells.Item(320, 7).Value = "TC4"
$ws.Cells.Item(320, 8).Value = "TCAO017"

$ws.Cells.Item(321, 1).Value = 2
$ws.Cells.Item(321, 2).Value = "Argyropelecus olfersii"
$ws.Cells.Item(321, 2).Font.Italic = $true
$ws.Cells.Item(321, 3).Value = "63mm SL"
$ws.Cells.Item(321, 4).Value = 7.43
$ws.Cells.Item(321, 5).Value = 2
$ws.Cells.Item(321, 6).Value = "B8"
$ws.Cells.Item(321, 7).Value = "TC4"
$ws.Cells.Item(321, 8).Value = "TCAO018"

$ws.Cells.Item(322, 1).Value = 2
$ws.Cells.Item(322, 2).Value = "Argyropelecus olfersii"
$ws.Cells.Item(322, 2).Font.Italic = $true
$ws.Cells.Item(322, 4).Value = 13.13
$ws.Cells.Item(322, 5).Value = 2
$ws.Cells.Item(322, 6).Value = "B9"
$ws.Cells.Item(322, 7).Value = "TC4"
$ws.Cells.Item(322, 8).Value = "TCAO019"
$ws.Cells.Item(322, 3).Value = "80mm SL"

$ws.Range("H322").Select()
$excel.ActiveWindow.ScrollRow = 308
$excel.ActiveWindow.ScrollColumn = 1
